$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.320.80'
$ws.Range("E2").Value = '  +1.26%  '

$ws.Range("D3").Value = '2.591.39'
$ws.Range("E3").Value = '  +4.69%  '

$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").Value = '''307.44'
$ws.Range("E5").Value = '  +4.67%  '

$ws.Range("D6").Value = '''100.29'
$ws.Range("E6").Value = '  +5.35%  '

$ws.Range("D7").Value = '''0.601'
$ws.Range("E7").Value = '  +4.76%  '

$ws.Range("E8").Value = '  +0.17%  '

$ws.Range("D9").Value = '''0.576'
$ws.Range("E9").Value = '  +10.30%  '

$ws.Range("D10").Value = '''39.32'
$ws.Range("E10").Value = '  +11.72%  '

$ws.Range("D11").Value = '''0.0845'
$ws.Range("E11").Value = '  +7.67%  '

$ws.Range("D12").Value = '''54.15'
$ws.Range("E12").Value = '  +0.82%  '

$ws.Range("D13").Value = '''8.15'
$ws.Range("E13").Value = '  +11.31%  '

$ws.Range("D14").Value = '3.003.75'
$ws.Range("E14").Value = '  +5.80%  '

$ws.Range("E15").Value = '  +1.83%  '

$ws.Range("D16").Value = '2.605.40'
$ws.Range("E16").Value = '  +5.18%  '

$ws.Range("D17").Value = '''0.919'
$ws.Range("E17").Value = '  +7.84%  '

$ws.Range("D18").Value = '''14.94'
$ws.Range("E18").Value = '  +4.95%  '

$ws.Range("D19").Value = '46.480.25'
$ws.Range("E19").Value = '  +1.38%  '

$ws.Range("E20").Value = '  +6.98%  '

$ws.Range("D21").Value = '''12.93'
$ws.Range("E21").Value = '  +1.44%  '

$ws.Range("D22").Value = '''6.71'
$ws.Range("E22").Value = '  +6.32%  '

$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '''71.44'
$ws.Range("E23").Value = '  +5.82%  '

$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").Value = '''274.80'
$ws.Range("E24").Value = '  +11.86%  '

$ws.Range("D25").Value = '''3.02'
$ws.Range("E25").Value = '  +8.32%  '

$ws.Range("E26").Value = '  +10.43%  '

$ws.Range("D27").Value = '''28.75'
$ws.Range("E27").Value = '  +31.31%  '

$ws.Range("D28").Value = '''1.00'
$ws.Range("E28").Value = '  -0.15%  '

$ws.Range("D29").Value = '''4.02'
$ws.Range("E29").Value = '  -0.71%  '

$ws.Range("D30").Value = '''10.58'
$ws.Range("E30").Value = '  +7.76%  '

$ws.Range("D31").Value = '''2.29'
$ws.Range("E31").Value = '  +3.57%  '

$ws.Range("D32").Value = '''38.90'
$ws.Range("E32").Value = '  -0.18%  '

$ws.Range("E33").Value = '  +14.57%  '

$ws.Range("D34").Value = '''3.62'
$ws.Range("E34").Value = '  -4.89%  '

$ws.Range("D35").Value = '''2.85'
$ws.Range("E35").Value = '  +3.00%  '

$ws.Range("D36").Value = '''2.21'
$ws.Range("E36").Value = '  +6.84%  '

$ws.Range("D37").Value = '''0.0838'
$ws.Range("E37").Value = '  +8.49%  '

$ws.Range("D38").Value = '''150.98'
$ws.Range("E38").Value = '  +2.31%  '

$ws.Range("E39").Value = '  +6.77%  '

$ws.Range("E40").Value = '  +5.96%  '

$ws.Range("D41").Value = '''23.36'
$ws.Range("E41").Value = '  +43.20%  '

$ws.Range("D42").Value = '''15.90'
$ws.Range("E42").Value = '  +4.37%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '''0.0332'
$ws.Range("E43").Value = '  +10.38%  '

$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").Value = '''3.64'
$ws.Range("E44").Value = '  +12.32%  '

$ws.Range("D45").Value = '''4.06'
$ws.Range("E45").Value = '  +2.66%  '

$ws.Range("D46").Value = '2.141.33'
$ws.Range("E46").Value = '  +6.26%  '

$ws.Range("D47").Value = '''0.997'
$ws.Range("E47").Value = '  -0.41%  '

$ws.Range("D48").Value = '''92.91'
$ws.Range("E48").Value = '  +1.55%  '

$ws.Range("D49").Value = '''9.47'
$ws.Range("E49").Value = '  +9.90%  '

$ws.Range("D50").Value = '''1.79'
$ws.Range("E50").Value = '  +1.47%  '

$ws.Range("D51").Value = '''108.58'
$ws.Range("E51").Value = '  +5.12%  '
